# Auto-generated edit script applying the Raiden_Profits.xlsx diff
# Mapping: ALC=sheet1, ARM=sheet2, BSM=sheet3, CRP=sheet4, CUL=sheet5, GSM=sheet6, LTW=sheet7, WVR=sheet8
$wb = $excel.ActiveWorkbook

# Row 19 (hunk 0)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1067.3889
$ws.Range("I19").Value = 1285.3334
$ws.Range("J19").Value = 849.44446
$ws.Range("K19").Value = 1285.3334
$ws.Range("L19").Value = 849.44446
$ws.Range("M19").Value = -1110.3334
$ws.Range("N19").Value = -1199.44446

# Row 43 (hunk 1)
$ws.Range("H43").Value = 18967.209
$ws.Range("I43").Value = 2238.2307
$ws.Range("K43").Value = 2238.2307
$ws.Range("M43").Value = -2169.2307

# Row 76 (hunk 2)
$ws.Range("H76").Value = 8996.4
$ws.Range("I76").Value = 8996.4
$ws.Range("K76").Value = 8996.4
$ws.Range("M76").Value = -8681.4

# Row 79 (hunk 3)
$ws.Range("H79").Value = 8996.4
$ws.Range("I79").Value = 8996.4
$ws.Range("K79").Value = 8996.4
$ws.Range("M79").Value = -7904.4

# Row 103 (hunk 4)
$ws.Range("H103").Value = 1491.3334
$ws.Range("J103").Value = 1491
$ws.Range("L103").Value = 4473
$ws.Range("N103").Value = -5645

# Row 132 (hunk 5)
$ws.Range("H132").Value = 627017.5
$ws.Range("I132").Value = 2308.76
$ws.Range("J132").Value = 2858120.2
$ws.Range("K132").Value = 6926.280000000001
$ws.Range("L132").Value = 8574360.600000001
$ws.Range("M132").Value = -4396.280000000001
$ws.Range("N132").Value = -8579420.600000001

# Row 138 (hunk 6)
$ws.Range("H138").Value = 2908.389
$ws.Range("J138").Value = 2979.7954
$ws.Range("L138").Value = 8939.386200000001
$ws.Range("N138").Value = -19219.3862

# Row 45 (hunk 7)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 911
$ws.Range("I45").Value = 911
$ws.Range("K45").Value = 911
$ws.Range("M45").Value = -534

# Row 122 (hunk 8)
$ws.Range("H122").Value = 1975.123
$ws.Range("I122").Value = 1831.4807
$ws.Range("J122").Value = 2549.6924
$ws.Range("K122").Value = 5494.4421
$ws.Range("L122").Value = 7649.0772
$ws.Range("M122").Value = -3044.4421
$ws.Range("N122").Value = -12549.0772

# Row 131 (hunk 9)
$ws.Range("H131").Value = 67500
$ws.Range("J131").Value = 67500
$ws.Range("L131").Value = 67500
$ws.Range("N131").Value = -77580

# Row 132 (hunk 10)
$ws.Range("H132").Value = 1360.1538
$ws.Range("J132").Value = 9998
$ws.Range("L132").Value = 29994
$ws.Range("N132").Value = -35054

# Row 25 (hunk 11)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 21199.8
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

# Row 134 (hunk 12)
$ws.Range("H134").Value = 4202.722
$ws.Range("I134").Value = 3143.2666
$ws.Range("J134").Value = 9500
$ws.Range("K134").Value = 9429.799800000001
$ws.Range("L134").Value = 28500
$ws.Range("M134").Value = -6894.799800000001
$ws.Range("N134").Value = -33570

# Row 140 (hunk 13)
$ws.Range("H140").Value = 78566.5
$ws.Range("J140").Value = 78566.5
$ws.Range("L140").Value = 78566.5
$ws.Range("N140").Value = -88926.5

# Row 141 (hunk 14)
$ws.Range("H141").Value = 95499.5
$ws.Range("J141").Value = 99333
$ws.Range("L141").Value = 99333
$ws.Range("N141").Value = -109693

# Row 16 (hunk 15)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 841.1667
$ws.Range("I16").Value = 784.5
$ws.Range("J16").Value = 1011.1667
$ws.Range("K16").Value = 784.5
$ws.Range("L16").Value = 1011.1667
$ws.Range("M16").Value = -497.5
$ws.Range("N16").Value = -1585.1667

# Row 113 (hunk 16)
$ws.Range("H113").Value = 841.1667
$ws.Range("I113").Value = 784.5
$ws.Range("J113").Value = 1011.1667
$ws.Range("K113").Value = 784.5
$ws.Range("L113").Value = 1011.1667
$ws.Range("M113").Value = 1385.5
$ws.Range("N113").Value = -5351.1667

# Row 134 (hunk 17)
$ws.Range("H134").Value = 1853.8889
$ws.Range("I134").Value = 1875.1333
$ws.Range("K134").Value = 5625.3999
$ws.Range("M134").Value = -3090.3999

# Row 75 (hunk 18)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 3647.6155
$ws.Range("I75").Value = 1894
$ws.Range("K75").Value = 5682
$ws.Range("M75").Value = -4684

# Row 78 (hunk 19)
$ws.Range("H78").Value = 3647.6155
$ws.Range("I78").Value = 1894
$ws.Range("K78").Value = 17046
$ws.Range("M78").Value = -12054

# Row 88 (hunk 20)
$ws.Range("H88").Value = 4184.6665
$ws.Range("J88").Value = 4927.5
$ws.Range("L88").Value = 14782.5
$ws.Range("N88").Value = -15638.5

# Row 91 (hunk 21)
$ws.Range("H91").Value = 4184.6665
$ws.Range("J91").Value = 4927.5
$ws.Range("L91").Value = 14782.5
$ws.Range("N91").Value = -17746.5

# Row 98 (hunk 22)
$ws.Range("H98").Value = 644.5
$ws.Range("J98").Value = 590
$ws.Range("L98").Value = 1770
$ws.Range("N98").Value = -4766

# Row 114 (hunk 23)
$ws.Range("H114").Value = 833.2273
$ws.Range("J114").Value = 868.8570999999999
$ws.Range("L114").Value = 2606.5713
$ws.Range("N114").Value = -9114.5713

# Row 125 (hunk 24)
$ws.Range("H125").Value = 9798.5
$ws.Range("I125").Value = 9796
$ws.Range("K125").Value = 29388
$ws.Range("M125").Value = -24468

# Row 131 (hunk 25)
$ws.Range("H131").Value = 3379426.8
$ws.Range("J131").Value = 5557740.5
$ws.Range("L131").Value = 16673221.5
$ws.Range("N131").Value = -16683301.5

# Row 80 (hunk 26)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 18606.777
$ws.Range("I80").Value = 2324.5
$ws.Range("K80").Value = 2324.5
$ws.Range("M80").Value = -1326.5

# Row 83 (hunk 27)
$ws.Range("H83").Value = 18606.777
$ws.Range("I83").Value = 2324.5
$ws.Range("K83").Value = 11622.5
$ws.Range("M83").Value = -6630.5

# Row 113 (hunk 28)
$ws.Range("H113").Value = 2257.1765
$ws.Range("I113").Value = 2124.8667
$ws.Range("K113").Value = 2124.8667
$ws.Range("M113").Value = 45.13329999999996

# Row 132 (hunk 29)
$ws.Range("H132").Value = 1635.4048
$ws.Range("I132").Value = 1544.421
$ws.Range("J132").Value = 2499.75
$ws.Range("K132").Value = 4633.263
$ws.Range("L132").Value = 7499.25
$ws.Range("M132").Value = -2103.263
$ws.Range("N132").Value = -12559.25

# Row 136 (hunk 30)
$ws.Range("H136").Value = 14007.23
$ws.Range("J136").Value = 14007.23
$ws.Range("L136").Value = 42021.69
$ws.Range("N136").Value = -47121.69

# Row 16 (hunk 31)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 905.1539
$ws.Range("I16").Value = 839
$ws.Range("K16").Value = 839
$ws.Range("M16").Value = -669

# Row 40 (hunk 32)
$ws.Range("H40").Value = 4633.5835
$ws.Range("I40").Value = 4043.5557
$ws.Range("J40").Value = 4987.6
$ws.Range("K40").Value = 4043.5557
$ws.Range("L40").Value = 4987.6
$ws.Range("M40").Value = -3907.5557
$ws.Range("N40").Value = -5259.6

# Row 55 (hunk 33)
$ws.Range("H55").Value = 382.81818
$ws.Range("I55").Value = 527.8570999999999
$ws.Range("K55").Value = 527.8570999999999
$ws.Range("M55").Value = -354.8570999999999

# Row 122 (hunk 34)
$ws.Range("H122").Value = 2710.375
$ws.Range("I122").Value = 2152.2222
$ws.Range("J122").Value = 3428
$ws.Range("K122").Value = 6456.6666
$ws.Range("L122").Value = 10284
$ws.Range("M122").Value = -4006.6666
$ws.Range("N122").Value = -15184

# Row 136 (hunk 35)
$ws.Range("H136").Value = 2179.4075
$ws.Range("I136").Value = 1567.7
$ws.Range("K136").Value = 4703.1
$ws.Range("M136").Value = -2153.1

# Row 30 (hunk 36)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 1000
$ws.Range("I30").Value = 1000
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 1000
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -893
$ws.Range("N30").ClearContents()

# Row 51 (hunk 37)
$ws.Range("H51").Value = 16000
$ws.Range("J51").Value = 16000
$ws.Range("L51").Value = 16000
$ws.Range("N51").Value = -17020

# Row 52 (hunk 38)
$ws.Range("H52").Value = 14509.167
$ws.Range("I52").Value = 7011
$ws.Range("J52").Value = 52000
$ws.Range("K52").Value = 7011
$ws.Range("L52").Value = 52000
$ws.Range("M52").Value = -6785
$ws.Range("N52").Value = -52452
